# Change "Ready for handoff" status text to "In Translation" across all sheets,
# then shrink the now-narrower Status columns to match the new auto-fit width.

$wb = $excel.ActiveWorkbook

# NOTE on $newWidth: the target OOXML column width is 13.4101845877511
# (an externally-computed, non-pixel-quantized auto-fit width). Excel's
# ColumnWidth COM setter snaps to the workbook's pixel grid (steps of 1/6
# of a character here), so that exact value is not reachable through the
# object model. 12.5 is the input that lands on the closest reachable
# grid point (13.333333333333334) to the target width.
$oldText = "Ready for handoff"
$newText = "In Translation"
$newWidth = 12.5

# Overview sheet: status text lives in the zh-cn (E) and de-de (F) columns, rows 2-3
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewStatusCells = @("E2", "F2", "E3", "F3")
foreach ($addr in $overviewStatusCells) {
    $cell = $wsOverview.Range($addr)
    if ($cell.Text -eq $oldText) {
        $cell.Value = $newText
    }
}

# zh-cn / de-de sheets: status text lives in the Status column (C), rows 2-3
foreach ($sheetName in @("zh-cn", "de-de")) {
    $sheet = $wb.Worksheets.Item($sheetName)
    foreach ($addr in @("C2", "C3")) {
        $cell = $sheet.Range($addr)
        if ($cell.Text -eq $oldText) {
            $cell.Value = $newText
        }
    }
}

# Overview sheet: zh-cn (E) and de-de (F) status columns
$wsOverview.Range("E1").EntireColumn.ColumnWidth = $newWidth
$wsOverview.Range("F1").EntireColumn.ColumnWidth = $newWidth

# zh-cn sheet: Status column (C)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = $newWidth

# de-de sheet: Status column (C)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = $newWidth
